# Shear Center Reference point added
# Update the Elements table so every section references the shear center
# reference point, then remove the now-unused rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-6 (columns A:D)
$data = @(
    @(0, 0, 1, 2),
    @(1, 1, 3, 1),
    @(2, 2, 3, 2),
    @(4, 3, 4, 1),
    @(5, 4, 5, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

# Remove the now-obsolete rows 7 through 10
$ws.Range("A7:D10").Delete()

# Update the active selection to reflect where the user left off
$ws.Range("F8").Select()
